$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.057.98"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "'2.300.75"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'300.30"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").Value = "'97.21"
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("D10").Value = "'33.65"
$ws.Range("E10").Value = "  -3.45%  "
$ws.Range("D11").Value = "'0.0794"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").Value = "'49.16"
$ws.Range("E12").Value = "  -4.06%  "
$ws.Range("E13").Value = "  +2.17%  "
$ws.Range("D14").Value = "'16.94"
$ws.Range("E14").Value = "  +9.60%  "
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").Value = "'2.658.31"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").Value = "'2.309.83"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "'0.810"
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("D19").Value = "'43.000.20"
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "'11.58"
$ws.Range("E21").Value = "  -1.05%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "'67.74"
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("D24").Value = "'236.71"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").Value = "'2.03"
$ws.Range("E25").Value = "  +3.41%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("D28").Value = "'24.48"
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("D29").Value = "'166.85"
$ws.Range("E29").Value = "  +1.28%  "
$ws.Range("D30").Value = "'33.96"
$ws.Range("E30").Value = "  -1.49%  "
$ws.Range("E31").Value = "  -5.81%  "
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "'4.69"
$ws.Range("E34").Value = "  +5.68%  "
$ws.Range("E35").Value = "  -1.32%  "
$ws.Range("E36").Value = "  -1.10%  "
$ws.Range("D37").Value = "'16.83"
$ws.Range("E37").Value = "  +2.84%  "
$ws.Range("D38").Value = "'0.0693"
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("D39").Value = "'2.83"
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("E41").Value = "  -1.73%  "
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("D43").Value = "'2.40"
$ws.Range("E43").Value = "  -1.54%  "
$ws.Range("D44").Value = "'1.980.32"
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("D46").Value = "'9.88"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").Value = "'17.68"
$ws.Range("E47").Value = "  -3.27%  "
$ws.Range("E48").Value = "  -0.87%  "
$ws.Range("D49").Value = "'2.525.80"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D50").Value = "'53.17"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("E51").Value = "  -4.18%  "
